# Apply schedule updates to Emerson's timetable (Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value  = "[-, -, 'MCT-3A-Robótica', -]"
$ws.Range("E2").Value  = "-"

$ws.Range("B3").Value  = "-"
$ws.Range("E3").Value  = "-"
$ws.Range("F3").Value  = "-"

$ws.Range("B4").Value  = "[-, -, 'MEC-2B-Elet. Dig. Bas.', -]"
$ws.Range("F4").Value  = "[-, 'MEC-2B-Elet. Dig. Bas.', -, 'MEC-1B-Comandos Eletricos']"

$ws.Range("B6").Value  = "['MCT-3A-Robótica', -, 'MEC-2B-Elet. Dig. Bas.', -]"
$ws.Range("F6").Value  = "[-, 'MEC-2B-Elet. Dig. Bas.', -, 'MEC-1B-Comandos Eletricos']"

$ws.Range("B7").Value  = "['MCT-3A-Robótica', -, -, -]"

$ws.Range("C8").Value  = "[-, -, -, 'MCT-3A-Robótica']"
$ws.Range("F8").Value  = "[-, -, -, 'MEC-1B-Comandos Eletricos']"

$ws.Range("D18").Value = "[-, 'ELM-2NA-Eletrônica Básica']"
$ws.Range("E18").Value = "[-, -, 'MEC-1NA-Comandos Eletricos', -]"

$ws.Range("D19").Value = "['ELM-2NA-Eletrônica Básica', 'ELM-2NA-Eletrônica Básica']"
$ws.Range("E19").Value = "[-, -, 'MEC-1NA-Comandos Eletricos', -]"
$ws.Range("F19").Value = "-"

$ws.Range("D20").Value = "['ELM-2NA-Eletrônica Básica', -]"
$ws.Range("E20").Value = "[-, -, 'MEC-1NA-Comandos Eletricos', -]"
$ws.Range("F20").Value = "-"

$ws.Range("E21").Value = "[-, -, 'MEC-1NA-Comandos Eletricos', -]"
$ws.Range("F21").Value = "-"
